$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.365.81"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.708.44"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.63"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5320"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2658"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06588"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.82"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07639"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.567"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("D13").Value = "1.719.88"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").Value = "1.945.39"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5721"
$ws.Range("E15").Value = "  -4.17%  "
$ws.Range("D16").Value = "0.0₅8160"
$ws.Range("E16").Value = "  -2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.82"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "27.373.85"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.21"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.667"
$ws.Range("E21").Value = "  -3.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("E22").Value = "  -4.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.969"
$ws.Range("E23").Value = "  -4.22%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.770"
$ws.Range("E25").Value = "  +6.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.87"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.271"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.30"
$ws.Range("E29").Value = "  -5.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05418"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.506"
$ws.Range("E32").Value = "  -5.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.426"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.644"
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.875"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9492"
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.412"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5859"
$ws.Range("E38").Value = "  -1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01632"
$ws.Range("E39").Value = "  -1.90%  "
$ws.Range("E40").Value = "  -1.20%  "
$ws.Range("D41").Value = "1.045.15"
$ws.Range("E41").Value = "  -0.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8431"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.96"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "1.852.12"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.97"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4497"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.059"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05243"
$ws.Range("E51").Value = "  -1.45%  "
